$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has two "PearsonLogo.png" inline pictures (one in the
# primary footer, one in the first-page footer) that are both currently
# named "image1.png" and need to become "image2.png", plus a
# "BTec_Logo-Orange" inline picture (in the first-page header) currently
# named "image2.jpg" that needs to become "image1.jpg".
#
# InlineShape has no settable Name property (matching real Word's object
# model), so each picture is round-tripped through ConvertToShape /
# ConvertToInlineShape to reach the settable Shape.Name property, which
# is what actually rewrites the <wp:docPr name="..."/> attribute.

function Rename-InlinePicture($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $null = $shape.ConvertToInlineShape()
}

# Primary footer -> footer2.xml (PearsonLogo, currently "image1.png")
$primaryFooterPic = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
Rename-InlinePicture $primaryFooterPic "image2.png"

# First-page footer -> footer1.xml (PearsonLogo, currently "image1.png")
$firstFooterPic = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
Rename-InlinePicture $firstFooterPic "image2.png"

# First-page header -> header1.xml (BTec_Logo-Orange, currently "image2.jpg")
$firstHeaderPic = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
Rename-InlinePicture $firstHeaderPic "image1.jpg"
